$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.93
$ws.Range("AB5").Value = 350
$ws.Range("AH5").Value = 23
$ws.Range("AI5").Value = 11.5
$ws.Range("AK5").Value = 8
$ws.Range("AT5").Value = 7.4
$ws.Range("AV5").Value = 70
$ws.Range("AX5").Value = 4.05
$ws.Range("BB5").Value = 17.5
$ws.Range("BC5").Value = 100

# Row 6
$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 3.7
$ws.Range("M6").Value = 1.14
$ws.Range("N6").Value = 5.5
$ws.Range("S6").Value = 1.67
$ws.Range("T6").Value = 2.1
$ws.Range("U6").Value = 2.38
$ws.Range("V6").Value = 1.53
$ws.Range("W6").Value = 5
$ws.Range("X6").Value = 9
$ws.Range("Z6").Value = 21
$ws.Range("AI6").Value = 17
$ws.Range("AT6").Value = 2.1
$ws.Range("AX6").Value = 5.5

# Row 7
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.67
$ws.Range("AW7").Value = 151

# Row 8
$ws.Range("O8").Value = 1.44
$ws.Range("P8").Value = 2.63
$ws.Range("Q8").Value = 2.4
$ws.Range("R8").Value = 1.53

# Row 15
$ws.Range("G15").Value = 1.48
$ws.Range("I15").Value = 7
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 11
$ws.Range("Q15").Value = 1.8
$ws.Range("R15").Value = 2
$ws.Range("W15").Value = 7
$ws.Range("AA15").Value = 12
$ws.Range("AB15").Value = 26
$ws.Range("AC15").Value = 11
$ws.Range("AD15").Value = 8
$ws.Range("AE15").Value = 17
$ws.Range("AG15").Value = 301
$ws.Range("AH15").Value = 17
$ws.Range("AI15").Value = 34
$ws.Range("AO15").Value = 7.5
$ws.Range("AV15").Value = 51
$ws.Range("BA15").Value = 126

# Row 17
$ws.Range("G17").Value = 2.05
$ws.Range("H17").Value = 3.25
$ws.Range("I17").Value = 3.1
$ws.Range("J17").Value = 2.88
$ws.Range("L17").Value = 3.75
$ws.Range("M17").Value = 1.05
$ws.Range("N17").Value = 11
$ws.Range("O17").Value = 1.29
$ws.Range("P17").Value = 3.5
$ws.Range("R17").Value = 1.85
$ws.Range("AL17").Value = 26
$ws.Range("AN17").Value = 4.33
$ws.Range("AX17").Value = 5
$ws.Range("AZ17").Value = 26
$ws.Range("BA17").Value = 51

# Row 46
$ws.Range("G46").Value = 2.35
$ws.Range("I46").Value = 2.67
$ws.Range("J46").Value = 2.92
$ws.Range("K46").Value = 2.15
$ws.Range("L46").Value = 3.2
$ws.Range("M46").Value = 1.02
$ws.Range("N46").Value = 12
$ws.Range("O46").Value = 1.21
$ws.Range("P46").Value = 3.55
$ws.Range("U46").Value = 1.53
$ws.Range("V46").Value = 2.2
$ws.Range("X46").Value = 13
$ws.Range("Y46").Value = 9.25
$ws.Range("Z46").Value = 25
$ws.Range("AA46").Value = 17.5
$ws.Range("AB46").Value = 23
$ws.Range("AE46").Value = 12
$ws.Range("AF46").Value = 45
$ws.Range("AH46").Value = 11
$ws.Range("AI46").Value = 15.5
$ws.Range("AJ46").Value = 10
$ws.Range("AK46").Value = 32
$ws.Range("AL46").Value = 20
$ws.Range("AM46").Value = 24
$ws.Range("AN46").Value = 4.4
$ws.Range("AO46").Value = 12
$ws.Range("AP46").Value = 19
$ws.Range("AQ46").Value = 50
$ws.Range("AX46").Value = 4.7
$ws.Range("AY46").Value = 14
$ws.Range("AZ46").Value = 19.5
$ws.Range("BA46").Value = 60
$ws.Range("BB46").Value = 80
$ws.Range("BC46").Value = 200

# Row 61
$ws.Range("N61").Value = 8
$ws.Range("O61").Value = 1.44
$ws.Range("P61").Value = 2.63
$ws.Range("Q61").Value = 2.35
$ws.Range("R61").Value = 1.57

# Row 62
$ws.Range("Q62").Value = 2.5
$ws.Range("R62").Value = 1.5

# Row 63
$ws.Range("G63").Value = 2.3
$ws.Range("M63").Value = 1.06
$ws.Range("N63").Value = 10
$ws.Range("Q63").Value = 2.15
$ws.Range("R63").Value = 1.67
$ws.Range("AN63").Value = 4.33
$ws.Range("AO63").Value = 13
